$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + data in column B
$ws.Range("B1").Value = "Resultado"
$ws.Range("B2").Value = "Anna’s Archive"

# Column widths (A grew from ~10.57 to 19; B is new at ~18.43)
$ws.Columns.Item(1).ColumnWidth = 18.166666666666668
$ws.Columns.Item(2).ColumnWidth = 17.592447916666668

# Selection moves to A14, with A2:A14 selected
$ws.Range("A2:A14").Select() | Out-Null
